# The workbook tracks daily "Pepino ensalada" price observations for the
# "Feria Lagunitas de Puerto Montt" market. A new daily observation is
# inserted as row 72 (pushing every following row down by one), so the
# data stays sorted the way the source feed expects.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 72, shifting rows 72:157
# down to 73:158 (and carrying their formatting with them, as Excel does).
$ws.Rows.Item(72).Insert()

# Populate the newly inserted row 72 with the new observation.
$ws.Range("A72").Value = 4
$ws.Range("B72").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C72").Value = "Los Lagos"
$ws.Range("D72").Value = 44483
$ws.Range("E72").Value = 10
$ws.Range("F72").Value = 100112043
$ws.Range("G72").Value = "Pepino ensalada"
$ws.Range("H72").Value = "Sin especificar"
$ws.Range("I72").Value = "Primera"
$ws.Range("J72").Value = 250
$ws.Range("K72").Value = 17000
$ws.Range("L72").Value = 17000
$ws.Range("M72").Value = 17000
$ws.Range("N72").Value = "`$/caja 60 unidades"
$ws.Range("O72").Value = "Región de Arica y Parinacota"
$ws.Range("P72").Value = 283
$ws.Range("Q72").Value = 60
$ws.Range("R72").Value = "Hortaliza"
